$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update headers
$ws.Range("A1").Value = "Week"
$ws.Range("B1").Value = "Predicted_Quantity"

# Week labels and predicted quantities
$weeks = @("2025-W43", "2025-W44", "2025-W45", "2025-W46", "2025-W47", "2025-W48", "2025-W49", "2025-W50")
$values = @(3, 2, 2, 2, 1, 1, 1, 1)

for ($i = 0; $i -lt $weeks.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $weeks[$i]
    $ws.Cells.Item($row, 2).Value = $values[$i]
}

# Clear column C entirely (was Predicted_Quantity, now removed)
$ws.Range("C1:C9").Clear()
